# Update the quadratic-svm-score prediction values (column B) with the
# freshly computed scores from the latest run of outputs-r202 /
# ful-path.csv (replacing the previous placeholder value of 1 for every
# row). Column A / column C are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.6934348384994156
$ws.Range("B3").Value = 0.85624048232744698
$ws.Range("B4").Value = 1.4632562427414726
$ws.Range("B5").Value = 1.7487631742460241
$ws.Range("B6").Value = 1.6934348384994156
$ws.Range("B7").Value = 0.85624048232744698
$ws.Range("B8").Value = 1.6848216597237897
$ws.Range("B9").Value = 1.3142429882819044
